$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new package entry row (row 5), mirroring row 4 but with Type = "ContentPackage"
$ws.Range("A5").Value = "GZIP"
$ws.Range("B5").Value = "GZIP"
$ws.Range("C5").Value = "GZIP"
$ws.Range("D5").Value = "1.0.0"
$ws.Range("E5").Value = "ContentPackage"

# F5 mirrors F4's date-string value exactly (same text, same default style) -
# copy F4 so Excel doesn't reinterpret the text as a date serial number.
$ws.Range("F4").Copy()
$ws.Range("F5").PasteSpecial(-4163)  # xlPasteValues
